# Adds the "Assumptions" body text, and new "Dependencies" / "Constraints"
# sections to the risks document, right before the trailing (bookmarked)
# empty paragraph - matching the target diff.
#
# The last paragraph in the document is an otherwise-empty paragraph that
# carries the _GoBack bookmark. Typing into that paragraph inserts content
# right before the bookmark, and splitting it with InsertParagraphAfter()
# pushes the (still-empty) bookmark paragraph further down while leaving
# the typed text behind in its own, new paragraph - exactly like a user
# placing their cursor at that point in Word and typing normally. We keep
# repeating that "type, then split" step, which means the bookmark
# paragraph always remains the very last paragraph in the document.

$d = $word.ActiveDocument

function Get-TypingRange($doc) {
    # The paragraph we want to keep typing into is always the one
    # immediately before the (still empty) trailing bookmark paragraph.
    $idx = $doc.Paragraphs.Count - 1
    $rng = $doc.Paragraphs.Item($idx).Range.Duplicate
    $rng.Collapse(0)
    return $rng
}

function Add-PlainRun($range, [string]$text) {
    $range.InsertAfter($text)
    $range.Font.Name = "Times New Roman"
    $range.Font.NameBi = "Times New Roman"
    $range.Collapse(0)
}

function Add-ItalicRedRun($range, [string]$text) {
    $range.InsertAfter($text)
    $range.Font.Name = "Times New Roman"
    $range.Font.NameBi = "Times New Roman"
    $range.Font.Italic = $true
    $range.Font.Color = 255
    $range.Collapse(0)
}

function Split-Paragraph($range) {
    $range.Collapse(0)
    $range.InsertParagraphAfter()
}

function Set-HeadingFormat($doc) {
    # Bold + underline the heading paragraph we just typed, including its
    # paragraph mark, so the formatting sticks on the pPr too.
    $idx = $doc.Paragraphs.Count - 1
    $p = $doc.Paragraphs.Item($idx).Range
    $p.Font.Name = "Times New Roman"
    $p.Font.NameBi = "Times New Roman"
    $p.Font.Bold = $true
    $p.Font.Underline = 1
}

# ---- Assumptions body text ----
$r = Get-TypingRange $d
Add-PlainRun $r "We assume that the user is able to use a simple and hopefully straightforward search interface. We also assume they are able to handle the resulting CSV ("
Add-ItalicRedRun $r "define"
Add-PlainRun $r ") files themselves."
Split-Paragraph $r

# ---- blank separator line ----
$r = Get-TypingRange $d
Split-Paragraph $r

# ---- "Dependencies" heading ----
$r = Get-TypingRange $d
Add-PlainRun $r "Dependencies"
Set-HeadingFormat $d
Split-Paragraph $r

# ---- Dependencies body text ----
$r = Get-TypingRange $d
Add-PlainRun $r "The only dependencies for this project are as follows. The users continue to classify the images and upload new ones, without this the client has nothing to gather data from. "
Add-PlainRun $r "The client has enough web hosting space to cope with an increasingly large data base."
Split-Paragraph $r

# ---- blank separator line ----
$r = Get-TypingRange $d
Split-Paragraph $r

# ---- "Constraints" heading ----
$r = Get-TypingRange $d
Add-PlainRun $r "Constraints"
Set-HeadingFormat $d
Split-Paragraph $r

# ---- Constraints body text ----
$r = Get-TypingRange $d
Add-PlainRun $r "The main constraint on this project is time. With a deadline in March there is only a short amount of time to get the system up and running. Due to this constraint some of the less core features may have to be set aside. So long as a working system is set up the project will have been successful. However some of the other features would greatly enhance the usability and user friendliness of the system so are well worth trying to implement if enough time is allocated."

# ---- trailing space inside the final (bookmarked) paragraph ----
$final = $d.Paragraphs.Last.Range
$final.Collapse(0)
Add-PlainRun $final " "
